# Scheduled-runner update: refresh computed market-profit figures across the Yojimbo_Profits job sheets.
# Mirrors the upstream commit "chore: update Sheets via scheduled runner":
# H/I/J/K/L/M/N price & profit columns are recomputed per-leve row; a couple of rows
# lose their HQ-profit (N) cell entirely once the HQ price drops to 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 93.59999999999999
$ws.Range("I8").Value = 93.59999999999999
$ws.Range("K8").Value = 280.8
$ws.Range("M8").Value = -141.8
$ws.Range("H64").Value = 3073.182
$ws.Range("I64").Value = 2780
$ws.Range("J64").Value = 3159.4119
$ws.Range("K64").Value = 2780
$ws.Range("L64").Value = 3159.4119
$ws.Range("M64").Value = -2532
$ws.Range("N64").Value = -3655.4119
$ws.Range("H67").Value = 3073.182
$ws.Range("I67").Value = 2780
$ws.Range("J67").Value = 3159.4119
$ws.Range("K67").Value = 2780
$ws.Range("L67").Value = 3159.4119
$ws.Range("M67").Value = -1922
$ws.Range("N67").Value = -4875.4119
$ws.Range("H76").Value = 3490.4814
$ws.Range("I76").Value = 3486.652
$ws.Range("K76").Value = 3486.652
$ws.Range("M76").Value = -3171.652
$ws.Range("H79").Value = 3490.4814
$ws.Range("I79").Value = 3486.652
$ws.Range("K79").Value = 3486.652
$ws.Range("M79").Value = -2394.652
$ws.Range("H100").Value = 2188.8
$ws.Range("I100").Value = 1486
$ws.Range("K100").Value = 1486
$ws.Range("M100").Value = -945
$ws.Range("H141").Value = 3610.652
$ws.Range("I141").Value = 3657.647
$ws.Range("J141").Value = 3477.5
$ws.Range("K141").Value = 10972.941
$ws.Range("L141").Value = 10432.5
$ws.Range("M141").Value = -5792.940999999999
$ws.Range("N141").Value = -20792.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1677.1428
$ws.Range("I45").Value = 1677.1428
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1677.1428
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1300.1428
$ws.Range("N45").Value = $null
$ws.Range("H63").Value = 2661.725
$ws.Range("I63").Value = 2537.7827
$ws.Range("J63").Value = 2829.4119
$ws.Range("K63").Value = 2537.7827
$ws.Range("L63").Value = 2829.4119
$ws.Range("M63").Value = -1851.7827
$ws.Range("N63").Value = -4201.4119
$ws.Range("H66").Value = 2661.725
$ws.Range("I66").Value = 2537.7827
$ws.Range("J66").Value = 2829.4119
$ws.Range("K66").Value = 12688.9135
$ws.Range("L66").Value = 14147.0595
$ws.Range("M66").Value = -9256.913500000001
$ws.Range("N66").Value = -21011.0595
$ws.Range("H88").Value = 3814.2856
$ws.Range("I88").Value = 3566.6667
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 3566.6667
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -3160.6667
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 3814.2856
$ws.Range("I91").Value = 3566.6667
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 3566.6667
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -2162.6667
$ws.Range("N91").Value = -6808
$ws.Range("H122").Value = 2126.4194
$ws.Range("I122").Value = 1813.2273
$ws.Range("J122").Value = 2892
$ws.Range("K122").Value = 5439.6819
$ws.Range("L122").Value = 8676
$ws.Range("M122").Value = -2989.6819
$ws.Range("N122").Value = -13576
$ws.Range("H132").Value = 2674.1345
$ws.Range("I132").Value = 2359.244
$ws.Range("K132").Value = 7077.732
$ws.Range("M132").Value = -4547.732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1666.5
$ws.Range("I105").Value = 1578.2354
$ws.Range("J105").Value = 2166.6667
$ws.Range("K105").Value = 1578.2354
$ws.Range("L105").Value = 2166.6667
$ws.Range("M105").Value = 168.7646
$ws.Range("N105").Value = -5660.6667
$ws.Range("H134").Value = 1582.3549
$ws.Range("I134").Value = 1414.4667
$ws.Range("J134").Value = 1739.75
$ws.Range("K134").Value = 4243.4001
$ws.Range("L134").Value = 5219.25
$ws.Range("M134").Value = -1708.4001
$ws.Range("N134").Value = -10289.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2400
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 2400
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H99").Value = 2091.5625
$ws.Range("I99").Value = 1484.6
$ws.Range("K99").Value = 1484.6
$ws.Range("M99").Value = 13.40000000000009
$ws.Range("H126").Value = 2091.5625
$ws.Range("I126").Value = 1484.6
$ws.Range("K126").Value = 4453.799999999999
$ws.Range("M126").Value = -1983.799999999999
$ws.Range("H132").Value = 2154.1
$ws.Range("I132").Value = 1498.9333
$ws.Range("K132").Value = 4496.7999
$ws.Range("M132").Value = -1966.7999
$ws.Range("H134").Value = 1552.7632
$ws.Range("I134").Value = 1586.4
$ws.Range("J134").Value = 1426.625
$ws.Range("K134").Value = 4759.200000000001
$ws.Range("L134").Value = 4279.875
$ws.Range("M134").Value = -2224.200000000001
$ws.Range("N134").Value = -9349.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62553.5
$ws.Range("J4").Value = 500000
$ws.Range("L4").Value = 1500000
$ws.Range("N4").Value = -1500224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4143.522
$ws.Range("I70").Value = 4138.2666
$ws.Range("J70").Value = 4153.375
$ws.Range("K70").Value = 4138.2666
$ws.Range("L70").Value = 4153.375
$ws.Range("M70").Value = -3868.2666
$ws.Range("N70").Value = -4693.375
$ws.Range("H73").Value = 4143.522
$ws.Range("I73").Value = 4138.2666
$ws.Range("J73").Value = 4153.375
$ws.Range("K73").Value = 4138.2666
$ws.Range("L73").Value = 4153.375
$ws.Range("M73").Value = -3202.2666
$ws.Range("N73").Value = -6025.375
$ws.Range("H80").Value = 2928
$ws.Range("I80").Value = 2934.6667
$ws.Range("J80").Value = 2923
$ws.Range("K80").Value = 2934.6667
$ws.Range("L80").Value = 2923
$ws.Range("M80").Value = -1936.6667
$ws.Range("N80").Value = -4919
$ws.Range("H83").Value = 2928
$ws.Range("I83").Value = 2934.6667
$ws.Range("J83").Value = 2923
$ws.Range("K83").Value = 14673.3335
$ws.Range("L83").Value = 14615
$ws.Range("M83").Value = -9681.333500000001
$ws.Range("N83").Value = -24599
$ws.Range("H122").Value = 2068.8667
$ws.Range("I122").Value = 1369.1
$ws.Range("J122").Value = 3468.4
$ws.Range("K122").Value = 4107.299999999999
$ws.Range("L122").Value = 10405.2
$ws.Range("M122").Value = -1657.299999999999
$ws.Range("N122").Value = -15305.2
$ws.Range("H126").Value = 1278
$ws.Range("I126").Value = 1263.3334
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 3790.0002
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1320.0002
$ws.Range("N126").Value = -8840
$ws.Range("H132").Value = 3134.1333
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 3223.5557
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 9670.667099999999
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -14730.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5899.2856
$ws.Range("I122").Value = 6572.5
$ws.Range("J122").Value = 5001.6665
$ws.Range("K122").Value = 19717.5
$ws.Range("L122").Value = 15004.9995
$ws.Range("M122").Value = -17267.5
$ws.Range("N122").Value = -19904.9995
$ws.Range("H132").Value = 4266.6665
$ws.Range("I132").Value = 3009.85
$ws.Range("J132").Value = 6200.231
$ws.Range("K132").Value = 9029.549999999999
$ws.Range("L132").Value = 18600.693
$ws.Range("M132").Value = -6499.549999999999
$ws.Range("N132").Value = -23660.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1378.3704
$ws.Range("I132").Value = 951.2
$ws.Range("K132").Value = 2853.6
$ws.Range("M132").Value = -323.6000000000004
